$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data (row 6)
$ws.Range("A6").Value = 135
$ws.Range("B6").Value = 35
$ws.Range("C6").Value = 1.5840000000000001
$ws.Range("D6").Value = 40.229999999999997
$ws.Range("E6").Value = 3.5950000000000002
$ws.Range("F6").Value = 91.32
$ws.Range("G6").Value = 1.4450000000000001
$ws.Range("H6").Value = 36.700000000000003
$ws.Range("I6").Value = 2.6080000000000001
$ws.Range("J6").Value = 66.239999999999995

# Update the selected cell to match the diff
$ws.Range("H7").Select()
